# This script replaces the paragraph "Cliente = static String cpf;" with new text
# ("Pet = static int contConsultas;") and inserts eleven new paragraphs after it,
# covering new exercises 30, 31 and 32 (per commit "mais dois exercicios de engenharia").
$d = $word.ActiveDocument

# Locate the paragraph to replace/extend
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Cliente = static String cpf;*") {
        $target = $p
        break
    }
}

if ($null -eq $target) {
    throw "Target paragraph 'Cliente = static String cpf;' not found"
}

# Build the WordprocessingML (flat-OPC) fragment with the replacement + new paragraphs,
# preserving exact run/paragraph formatting (color, size, italics, highlight, underline).
$paragraphsXml = @(
    '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="both"/><w:rPr><w:i w:val="false"/><w:iCs w:val="false"/><w:color w:val="003D73"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:i w:val="false"/><w:iCs w:val="false"/><w:color w:val="003D73"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Pet = static int contConsultas;</w:t></w:r></w:p>',
    '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="both"/><w:rPr><w:i w:val="false"/><w:iCs w:val="false"/><w:color w:val="003D73"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:i w:val="false"/><w:iCs w:val="false"/><w:color w:val="003D73"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:r></w:p>',
    '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="both"/><w:rPr><w:i w:val="false"/><w:iCs w:val="false"/><w:color w:val="434343"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:i w:val="false"/><w:iCs w:val="false"/><w:color w:val="434343"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>30- Apresente a estrutura básica de código em JAVA, C# ou C++ para implementar os seis membros estáticos.</w:t></w:r></w:p>',
    '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="both"/><w:rPr><w:i/><w:iCs/><w:color w:val="434343"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="434343"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Feito no Eclipse.</w:t></w:r></w:p>',
    '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="both"/><w:rPr><w:i/><w:iCs/><w:color w:val="434343"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="434343"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:r></w:p>',
    '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="both"/><w:rPr><w:i w:val="false"/><w:iCs w:val="false"/><w:color w:val="434343"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:i w:val="false"/><w:iCs w:val="false"/><w:color w:val="434343"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>31- Transforme todos os relacionamentos de associação ou agregação entre as classes de entidade e todos os relacionamentos entre as classes de fronteira e controle para dependências estruturais. Explique a vantagem e desvantagem desse tipo de dependência.</w:t></w:r></w:p>',
    '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="both"/><w:rPr><w:i/><w:iCs/><w:color w:val="434343"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="434343"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>A dependência estrutural é um tipo de dependência mais fácil de ser implementado e pensado. Nessa dependencia, a classe dependente possui um atribuo que referencia outra classe.</w:t></w:r></w:p>',
    '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="both"/><w:rPr><w:i/><w:iCs/><w:color w:val="434343"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="434343"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Isso acaba ocasionando um alto acoplamento (baixa coesão). Por isso, precisamos observar se as dependências das classes não podem ser modificadas para depedências não estruturais.</w:t></w:r></w:p>',
    '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="both"/><w:rPr><w:i/><w:iCs/><w:color w:val="434343"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="434343"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="yellow"/></w:rPr><w:t>FALTA COLOCAR AS CLASSES DE FRONTEIRA NO DIAGRAMA E DEPOIS TERMINAR ESSE EXERCICIO</w:t></w:r></w:p>',
    '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="both"/><w:rPr><w:i/><w:iCs/><w:color w:val="434343"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="434343"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:r></w:p>',
    '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="both"/><w:rPr><w:i w:val="false"/><w:iCs w:val="false"/><w:color w:val="434343"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:i w:val="false"/><w:iCs w:val="false"/><w:color w:val="434343"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="none"/></w:rPr><w:t>32- Apresente a estrutura básica de código em JAVA, C# ou C++ para implementar as dependências estruturais.</w:t></w:r></w:p>',
    '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="both"/><w:rPr><w:i/><w:iCs/><w:color w:val="434343"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="yellow"/><w:u w:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="434343"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="yellow"/><w:u w:val="none"/></w:rPr><w:t>Já foi feito nos exercícios anteriores… acredito que tenhamos que focar apenas em melhorar nosso diagrama...</w:t></w:r></w:p>'
) -join ""

$xml = "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" +
    "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" +
    "<pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">" +
    "<w:body>" + $paragraphsXml + "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

# InsertXML REPLACES the target range's contents with the supplied WordOpenXML,
# so this single call both rewrites the first paragraph and appends all new ones.
$target.Range.InsertXML($xml)

Write-Host "Replacement complete"
